# Updated cryptos list on Sat Jan  6 20:50:15 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value while preserving/forcing plain-text storage
# (avoids numeric-looking strings like "94.66" being auto-coerced into
# real numbers by the Value setter), without leaving behind a lingering
# custom cell style on the written cells.
function Set-TextValue($ws, $addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue $ws "D2" "44.191.22"
Set-TextValue $ws "E2" "  +0.38%  "

# Row 3 - Ethereum
Set-TextValue $ws "D3" "2.242.21"
Set-TextValue $ws "E3" "  +0.04%  "

# Row 4 - TetherUSD
Set-TextValue $ws "E4" "  +0.13%  "

# Row 5 - BNB
Set-TextValue $ws "D5" "306.88"
Set-TextValue $ws "E5" "  -2.59%  "

# Row 6 - Solana
Set-TextValue $ws "D6" "94.66"
Set-TextValue $ws "E6" "  -4.61%  "

# Row 7 - XRP
Set-TextValue $ws "D7" "0.571"
Set-TextValue $ws "E7" "  -0.46%  "

# Row 8 - USDC
Set-TextValue $ws "D8" "1.01"
Set-TextValue $ws "E8" "  +0.29%  "

# Row 9 - Cardano
Set-TextValue $ws "D9" "0.525"
Set-TextValue $ws "E9" "  -1.44%  "

# Row 10 - Avalanche
Set-TextValue $ws "D10" "34.69"
Set-TextValue $ws "E10" "  -4.18%  "

# Row 11 - Dogecoin
Set-TextValue $ws "D11" "0.0810"
Set-TextValue $ws "E11" "  -1.59%  "

# Row 12 - Polkadot
Set-TextValue $ws "D12" "7.18"
Set-TextValue $ws "E12" "  -2.67%  "

# Row 13 - TRON
Set-TextValue $ws "E13" "  -0.21%  "

# Rows 14 & 15 swap places: WrappedEther <-> WrappedliquidstakedEther2.0
Set-TextValue $ws "B14" "WrappedliquidstakedEther2.0"
Set-TextValue $ws "C14" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws "D14" "2.585.97"
Set-TextValue $ws "E14" "  +0.15%  "

Set-TextValue $ws "B15" "WrappedEther"
Set-TextValue $ws "C15" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws "D15" "2.327.07"
Set-TextValue $ws "E15" "  +3.58%  "

# Row 16 - Polygon
Set-TextValue $ws "D16" "0.830"
Set-TextValue $ws "E16" "  -1.52%  "

# Row 17 - Chainlink
Set-TextValue $ws "D17" "13.52"
Set-TextValue $ws "E17" "  -3.09%  "

# Row 18 - WrappedBTC
Set-TextValue $ws "D18" "43.976.21"
Set-TextValue $ws "E18" "  +0.24%  "

# Row 19 - ShibaInu
Set-TextValue $ws "D19" "0.0₃0962"
Set-TextValue $ws "E19" "  -1.27%  "

# Row 20 - Uniswap
Set-TextValue $ws "D20" "6.37"
Set-TextValue $ws "E20" "  +0.88%  "

# Row 21 - InternetComputer(DFINITY)
Set-TextValue $ws "D21" "12.14"
Set-TextValue $ws "E21" "  -8.10%  "

# Row 22 - Litecoin
Set-TextValue $ws "D22" "65.53"
Set-TextValue $ws "E22" "  -0.69%  "

# Row 23 - BitcoinCash
Set-TextValue $ws "D23" "237.83"
Set-TextValue $ws "E23" "  +0.04%  "

# Row 24 - PancakeSwap
Set-TextValue $ws "D24" "2.95"
Set-TextValue $ws "E24" "  -0.61%  "

# Row 25 - ImmutableX
Set-TextValue $ws "D25" "2.00"
Set-TextValue $ws "E25" "  -1.10%  "

# Row 26 - Dai
Set-TextValue $ws "E26" "  +0.09%  "

# Row 27 - Cosmos
Set-TextValue $ws "D27" "9.93"
Set-TextValue $ws "E27" "  -2.16%  "

# Row 28 - InjectiveProtocol
Set-TextValue $ws "D28" "38.22"
Set-TextValue $ws "E28" "  +4.63%  "

# Row 29 - Toncoin
Set-TextValue $ws "E29" "  +3.17%  "

# Row 30 - EthereumClassic
Set-TextValue $ws "D30" "20.04"
Set-TextValue $ws "E30" "  -0.06%  "

# Row 31 - Filecoin
Set-TextValue $ws "D31" "5.85"
Set-TextValue $ws "E31" "  -2.36%  "

# Row 32 - Monero
Set-TextValue $ws "D32" "153.24"
Set-TextValue $ws "E32" "  -0.79%  "

# Row 33 - Hedera
Set-TextValue $ws "D33" "0.0795"
Set-TextValue $ws "E33" "  -4.78%  "

# Row 34 - WEMIXToken
Set-TextValue $ws "E34" "  -1.71%  "

# Row 35 - LidoDAOToken
Set-TextValue $ws "D35" "3.17"
Set-TextValue $ws "E35" "  -5.00%  "

# Row 36 - Stellar
Set-TextValue $ws "E36" "  +1.74%  "

# Row 37 - Kaspa
Set-TextValue $ws "E37" "  -0.76%  "

# Row 38 - ARBITRUM
Set-TextValue $ws "D38" "1.76"
Set-TextValue $ws "E38" "  -7.78%  "

# Row 39 - NEARProtocol
Set-TextValue $ws "D39" "3.51"
Set-TextValue $ws "E39" "  -0.33%  "

# Row 40 - RenderToken
Set-TextValue $ws "D40" "3.82"
Set-TextValue $ws "E40" "  -4.74%  "

# Row 41 - Celestia
Set-TextValue $ws "D41" "14.41"
Set-TextValue $ws "E41" "  -8.71%  "

# Row 42 - VeChain
Set-TextValue $ws "D42" "0.0299"
Set-TextValue $ws "E42" "  -2.90%  "

# Row 43 - FirstDigitalUSD (unchanged)

# Row 44 - Maker
Set-TextValue $ws "D44" "1.746.47"
Set-TextValue $ws "E44" "  +2.37%  "

# Row 45 - BitcoinSV
Set-TextValue $ws "D45" "82.51"
Set-TextValue $ws "E45" "  +0.46%  "

# Row 46 - Algorand
Set-TextValue $ws "D46" "0.191"
Set-TextValue $ws "E46" "  -2.36%  "

# Row 47 - Aave
Set-TextValue $ws "D47" "99.74"
Set-TextValue $ws "E47" "  -2.00%  "

# Row 48 - THORChain
Set-TextValue $ws "D48" "4.93"
Set-TextValue $ws "E48" "  -4.69%  "

# Row 49 - FraxShare
Set-TextValue $ws "D49" "8.07"
Set-TextValue $ws "E49" "  -0.87%  "

# Row 50 - Stacks
Set-TextValue $ws "D50" "1.57"
Set-TextValue $ws "E50" "  -3.10%  "

# Row 51 - MultiversX
Set-TextValue $ws "D51" "54.42"
Set-TextValue $ws "E51" "  -3.52%  "
